$wb = $excel.ActiveWorkbook

# Add the new worksheet "ODI Batting Extra" after the last existing sheet
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "ODI Batting Extra"

# Match the look & feel (outline/page setup) used by the rest of the workbook
$ws.Outline.SummaryRow = 1
$ws.Outline.SummaryColumn = 1
$ps = $ws.PageSetup
$ps.LeftMargin = 54
$ps.RightMargin = 54
$ps.TopMargin = 72
$ps.BottomMargin = 72
$ps.HeaderMargin = 36
$ps.FooterMargin = 36

# Headers (row 1)
$ws.Range("A1").Value = "MATCH_CODE"
$ws.Range("B1").Value = "BATTING_POSITION"
$ws.Range("C1").Value = "NUM_4"
$ws.Range("D1").Value = "NUM_6"
$ws.Range("E1").Value = "PERCENT_RUNS_OF_TOTAL"
$ws.Range("F1").Value = "MAN_OF_MATCH"

# Re-use the header style already used on the other sheets (bold, centered, bordered)
$wb.Worksheets.Item("ODI Bowling").Range("A1").Copy()
$ws.Range("A1:F1").PasteSpecial(-4122)  # xlPasteFormats

# Data (row 2) -- keep as text where the source data is textual
$ws.Range("A2").NumberFormat = "@"
$ws.Range("C2").NumberFormat = "@"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("F2").NumberFormat = "@"

$ws.Range("A2").Value = "4717"
$ws.Range("B2").Value = 10
$ws.Range("C2").Value = "0"
$ws.Range("D2").Value = "0"
$ws.Range("E2").Value = "1.02%"
$ws.Range("F2").Value = "NO"

# Remove the temporary number-format styling so the cells fall back to the
# default (unstyled) cell, matching the source data which carries no style.
$ws.Range("A2").ClearFormats()
$ws.Range("C2").ClearFormats()
$ws.Range("D2").ClearFormats()
$ws.Range("E2").ClearFormats()
$ws.Range("F2").ClearFormats()

# Keep the original active sheet/selection ("Player Info") as it was before edit
$ws.Range("A1").Select() | Out-Null
$wb.Worksheets.Item("Player Info").Activate() | Out-Null

$wb.Save()
